# Apply the numeric updates described in the diff.
# Each old multiplication string is unique in the document, so a simple
# Find/Replace (wildcards off, match whole string) for each pair is safe.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "75×37=2775"; New = "31×33=1023" },
    @{ Old = "60×81=4860"; New = "67×39=2613" },
    @{ Old = "38×66=2508"; New = "67×91=6097" },
    @{ Old = "61×95=5795"; New = "43×99=4257" },
    @{ Old = "46×40=1840"; New = "56×40=2240" },
    @{ Old = "69×17=1173"; New = "17×61=1037" },
    @{ Old = "68×13=884";  New = "68×82=5576" },
    @{ Old = "55×51=2805"; New = "79×18=1422" },
    @{ Old = "39×61=2379"; New = "60×30=1800" },
    @{ Old = "63×33=2079"; New = "82×51=4182" },
    @{ Old = "96×37=3552"; New = "59×44=2596" },
    @{ Old = "24×99=2376"; New = "48×72=3456" },
    @{ Old = "91×36=3276"; New = "32×97=3104" },
    @{ Old = "88×90=7920"; New = "81×79=6399" },
    @{ Old = "48×50=2400"; New = "84×38=3192" },
    @{ Old = "90×50=4500"; New = "97×53=5141" },
    @{ Old = "33×35=1155"; New = "70×26=1820" },
    @{ Old = "43×71=3053"; New = "36×24=864" },
    @{ Old = "11×12=132";  New = "87×45=3915" },
    @{ Old = "61×31=1891"; New = "48×84=4032" },
    @{ Old = "34×63=2142"; New = "91×41=3731" },
    @{ Old = "50×52=2600"; New = "45×27=1215" },
    @{ Old = "51×65=3315"; New = "49×33=1617" },
    @{ Old = "20×49=980";  New = "14×26=364" },
    @{ Old = "43×65=2795"; New = "18×86=1548" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
